$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 799, shifting the existing rows 799:840 down to 800:841
$ws.Rows.Item(799).Insert()

# Populate the newly inserted row with its data.
# Column A holds a date formatted as plain text (e.g. "2026/02/08"); format the
# cell as Text first so Excel's COM layer doesn't auto-convert the literal
# into a real date serial number, then restore the default "Normal" style so
# no stray number-format styling is left behind on the cell.
$ws.Cells.Item(799, 1).NumberFormat = "@"
$ws.Cells.Item(799, 1).Value = "2026/02/08"
$ws.Cells.Item(799, 1).Style = "Normal"

$ws.Cells.Item(799, 2).Value = "日"
$ws.Cells.Item(799, 3).Value = 17
$ws.Cells.Item(799, 4).Value = 201
